$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Update the "last updated" timestamp banner (A1) ---
$ws.Range("A1").Value = "Datos actualizados a 29 de Julio de 2020 a las 23:22"

# --- Country rank swaps: two countries' live data crossed over, so the
#     two (or three) rows exchange which country label they show, while
#     new refreshed numbers land on the row that now sits higher. ---

# Etiopia / Costa de Marfil (rows 74-75)
$ws.Range("A74").Value = "Costa de Marfil"
$ws.Range("A75").Value = "Etiopia"

# Guinea-Bisau / Ruanda (rows 129-130)
$ws.Range("A129").Value = "Ruanda"
$ws.Range("A130").Value = "Guinea-Bisau"

# Camboya / Islas Feroe / Guadalupe (rows 175-177) - Guadalupe jumped two spots
$ws.Range("A175").Value = "Guadalupe"
$ws.Range("A176").Value = "Camboya"
$ws.Range("A177").Value = "Islas Feroe"

# --- Refreshed numeric stats (Casos totales, Nuevos casos, Casos activos,
#     Recuperados, Casos criticos, Muertes hoy, Muertes) ---

# Row 4 - Estados Unidos
$ws.Range("B4").Value = 4549472
$ws.Range("C4").Value = 51129
$ws.Range("D4").Value = 2228978
$ws.Range("E4").Value = 2167131
$ws.Range("G4").Value = 1072
$ws.Range("H4").Value = 153363

# Row 5 - Brasil
$ws.Range("B5").Value = 2553265
$ws.Range("C5").Value = 68616
$ws.Range("E5").Value = 741571
$ws.Range("G5").Value = 1500
$ws.Range("H5").Value = 90134

# Row 8 - Sudafrica
$ws.Range("B8").Value = 471123
$ws.Range("C8").Value = 11362
$ws.Range("D8").Value = 297967
$ws.Range("E8").Value = 165659
$ws.Range("G8").Value = 240
$ws.Range("H8").Value = 7497

# Row 25 - Canada
$ws.Range("B25").Value = 115298
$ws.Range("C25").Value = 304
$ws.Range("D25").Value = 100335
$ws.Range("E25").Value = 6049

# Row 28 - Egipto
$ws.Range("B28").Value = 93356
$ws.Range("C28").Value = 409
$ws.Range("D28").Value = 37025
$ws.Range("E28").Value = 51603
$ws.Range("G28").Value = 37
$ws.Range("H28").Value = 4728

# Row 36 - Israel
$ws.Range("B36").Value = 68299
$ws.Range("C36").Value = 2006
$ws.Range("D36").Value = 32746
$ws.Range("E36").Value = 35062

# Row 52 - Barein
$ws.Range("B52").Value = 40311
$ws.Range("C52").Value = 390
$ws.Range("D52").Value = 36920
$ws.Range("E52").Value = 3247
$ws.Range("G52").Value = 3
$ws.Range("H52").Value = 144

# Row 55 - Suiza
$ws.Range("D55").Value = 31100
$ws.Range("E55").Value = 1723

# Row 74 - now Costa de Marfil (fresh data)
$ws.Range("B74").Value = 15813
$ws.Range("C74").Value = 100
$ws.Range("D74").Value = 10793
$ws.Range("E74").Value = 4921
$ws.Range("G74").Value = 1
$ws.Range("H74").Value = 99

# Row 75 - now Etiopia (carries forward old Etiopia numbers)
$ws.Range("B75").Value = 15810
$ws.Range("C75").Value = 610
$ws.Range("D75").Value = 6685
$ws.Range("E75").Value = 8872
$ws.Range("G75").Value = 14
$ws.Range("H75").Value = 253

# Row 80 - Estado de Palestina
$ws.Range("E80").Value = 6371
$ws.Range("G80").Value = 1
$ws.Range("H80").Value = 80

# Row 129 - now Ruanda (fresh data)
$ws.Range("B129").Value = 1963
$ws.Range("C129").Value = 37
$ws.Range("D129").Value = 1036
$ws.Range("E129").Value = 922
$ws.Range("H129").Value = 5

# Row 130 - now Guinea-Bisau (carries forward old Guinea-Bisau numbers)
$ws.Range("B130").Value = 1954
$ws.Range("D130").Value = 803
$ws.Range("E130").Value = 1125
$ws.Range("H130").Value = 26

# Row 135 - Yemen
$ws.Range("B135").Value = 1711
$ws.Range("C135").Value = 8
$ws.Range("D135").Value = 849
$ws.Range("E135").Value = 377
$ws.Range("G135").Value = 1
$ws.Range("H135").Value = 485

# Row 175 - now Guadalupe (fresh data)
$ws.Range("B175").Value = 244
$ws.Range("C175").Value = 41
$ws.Range("D175").Value = 176
$ws.Range("E175").Value = 54
$ws.Range("H175").Value = 14

# Row 176 - now Camboya (carries forward old Camboya numbers)
$ws.Range("B176").Value = 226
$ws.Range("D176").Value = 147
$ws.Range("E176").Value = 79

# Row 177 - now Islas Feroe (carries forward old Islas Feroe numbers)
$ws.Range("B177").Value = 220
$ws.Range("D177").Value = 188
$ws.Range("E177").Value = 32
$ws.Range("H177").Value = 0
